$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the MemberId column entirely (was column A) - everything shifts left.
$ws.Columns.Item(1).Delete()

# After the shift the columns are:
#   A = FullName, B = Address, C = MovieId (-> MoviesRented), D = SalutationId (-> Salutation)

# Rename the headers.
$ws.Cells.Item(1, 3).Value = "MoviesRented"
$ws.Cells.Item(1, 4).Value = "Salutation"

# Replace the numeric MovieId values with the rented movie titles.
$ws.Cells.Item(2, 3).Value = "Daddy's Little Girls"
$ws.Cells.Item(3, 3).Value = "Clash of the Titans 2"
$ws.Cells.Item(4, 3).Value = "Forgetting Sarah Marshal"
$ws.Cells.Item(5, 3).Value = "Clash of the Titans 2"
$ws.Cells.Item(6, 3).Value = "Daddy's Little Girls"

# Replace the numeric SalutationId values with the salutation text.
$ws.Cells.Item(2, 4).Value = "Ms"
$ws.Cells.Item(3, 4).Value = "Ms"
$ws.Cells.Item(4, 4).Value = "Mr"
$ws.Cells.Item(5, 4).Value = "Mr"
$ws.Cells.Item(6, 4).Value = "Mr"
